$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new shared string / cell values for row 40
$ws.Range("A40").Value = 138
$ws.Range("B40").Value = "Copy List with Random Pointer"
$ws.Range("C40").Value = "Java"

# Match formatting of the row above (row 39): A has style index 10 (centered/top, green fill),
# B and C have style index 12 (green fill only)
$ws.Range("A39").Copy()
$ws.Range("A40").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B39:C39").Copy()
$ws.Range("B40:C40").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# Update selection to mirror the recorded state after edit
$ws.Range("B47").Select()
